$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add a new user row (row 6) below the existing data: a new user/email pair
# that reuses the same password text as row 5 ("1234").
$ws.Range("A6").Value = "navichoque"
$ws.Range("B6").Value = "dnavichoque@gmail.com"

# Copy the password cell from row 5 so the new cell keeps the same text
# (shared-string) type and style as the rest of the table instead of being
# re-interpreted as a number.
$ws.Range("C5").Copy($ws.Range("C6"))
